# TC45_Canine_Filter_Breed-Weimrnr.xlsx -- "startup" sheet query fix-up.
#
# The CasesTab query (cell B2) previously returned an extra `Cohort` column
# sourced from an OPTIONAL MATCH on (co:cohort) that isn't declared in the
# query's own MATCH clauses any more -- drop the trailing
# ", coalesce(co.cohort_description, '') AS `Cohort`" fragment so the
# RETURN clause ends cleanly after `Response to Treatment`.
#
# The FilesTab query (cell B4) text itself is unchanged; rewriting it here
# (with its own unchanged text) simply mirrors the workbook author's
# save, which also reshuffled/deduplicated the shared-string table so the
# two query strings end up swapping slots.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 'MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed IN [''Weimaraner'']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '''') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '''') AS `Study Code` ,
        coalesce(s.clinical_study_type, '''') AS  `Study Type`,
        coalesce(demo.breed, '''') AS Breed ,
        coalesce(diag.disease_term, '''') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '''') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '''') AS Age ,
        coalesce(demo.sex, '''') AS Sex ,
        coalesce(demo.neutered_indicator, '''') AS `Neutered Status`,
        coalesce(demo.weight, '''') AS `Weight (kg)`,
        coalesce(diag.best_response, '''') AS `Response to Treatment`'

$ws.Range("B4").Value = 'MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN [''Weimaraner'']
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '''') AS `File Name`,
        coalesce(labels(parent)[0], '''') AS `Association`,
        coalesce(f.file_description, '''') AS `Description`,
        coalesce(f.file_format, '''') AS `Format`,
        coalesce(f.file_size, '''') AS `Size`,
        coalesce(c.case_id, '''') AS `Case ID`,
        coalesce(diag.disease_term,'''') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'''') AS `Study Code`'

# Wrapped-text rows reflow to the new (shorter) line counts for B2/B4.
$ws.Rows.Item(2).RowHeight = 244.8
$ws.Rows.Item(3).RowHeight = 216
$ws.Rows.Item(4).RowHeight = 216

# Active selection moves from the FilesTab row (B4) to the CasesTab row (B2).
$ws.Range("B2").Select()
